$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new columns before the existing "ExpPoints" column (C), which
# shifts ExpPoints from C to G and makes room for WIN / TOP4 / TOP5 / RELEGATION.
$ws.Range("C1:F1").EntireColumn.Insert()

# New header labels for the inserted columns.
$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP4"
$ws.Range("E1").Value = "TOP5"
$ws.Range("F1").Value = "RELEGATION"

# Tottenham Hotspur and AFC Bournemouth swapped rank positions (rows 8 & 9).
$ws.Range("B8").Value = "Tottenham Hotspur"
$ws.Range("B9").Value = "AFC Bournemouth"

# Updated ExpPoints values (now in column G) from the refreshed prediction run.
$ws.Range("G2").Value = 77.78090953290865
$ws.Range("G3").Value = 74.65130565425311
$ws.Range("G4").Value = 70.41620903482935
$ws.Range("G5").Value = 61.54979332109163
$ws.Range("G6").Value = 59.00039416365026
$ws.Range("G7").Value = 55.85494031174017
$ws.Range("G8").Value = 55.38199809117498
$ws.Range("G9").Value = 55.12724226780722
$ws.Range("G10").Value = 54.55690505803825
$ws.Range("G11").Value = 52.70319100480773
$ws.Range("G12").Value = 52.61966813208718
$ws.Range("G13").Value = 48.59662354545544
$ws.Range("G14").Value = 45.36842609598845
$ws.Range("G15").Value = 44.8999452509425
$ws.Range("G16").Value = 40.60850717270075
$ws.Range("G17").Value = 40.1091405840919
$ws.Range("G18").Value = 37.25284888681973
$ws.Range("G19").Value = 35.97022327669699
$ws.Range("G20").Value = 33.61382405199829
$ws.Range("G21").Value = 29.4880830383504
